# Auto-generated Excel COM-interop script
# Updates market-data columns (H-N) on multiple Leve-profit sheets
# to match the scheduled data refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1408.25
$ws.Range("I41").Value = 483
$ws.Range("J41").Value = 1716.6666
$ws.Range("K41").Value = 483
$ws.Range("L41").Value = 1716.6666
$ws.Range("M41").Value = -43
$ws.Range("N41").Value = -2596.6666

$ws.Range("H53").Value = 310.2
$ws.Range("I53").Value = 250.25
$ws.Range("K53").Value = 250.25
$ws.Range("M53").Value = 386.75

$ws.Range("H58").Value = 1641.25
$ws.Range("J58").Value = 2440
$ws.Range("L58").Value = 7320
$ws.Range("N58").Value = -7620

$ws.Range("H106").Value = 10005
$ws.Range("I106").Value = 10005
$ws.Range("K106").Value = 10005
$ws.Range("M106").Value = -9374

$ws.Range("H132").Value = 10269.818
$ws.Range("I132").Value = 11969.556
$ws.Range("J132").Value = 2621
$ws.Range("K132").Value = 35908.66800000001
$ws.Range("L132").Value = 7863
$ws.Range("M132").Value = -33378.66800000001
$ws.Range("N132").Value = -12923

$ws.Range("H140").Value = 40260
$ws.Range("J140").Value = 40260
$ws.Range("L140").Value = 40260
$ws.Range("N140").Value = -50620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4689.8
$ws.Range("I122").Value = 4485
$ws.Range("K122").Value = 13455
$ws.Range("M122").Value = -11005

$ws.Range("H126").Value = 5881.6665
$ws.Range("I126").Value = 5881.6665
$ws.Range("K126").Value = 17644.9995
$ws.Range("M126").Value = -15174.9995

$ws.Range("H132").Value = 3979.8
$ws.Range("I132").Value = 3812.5715
$ws.Range("J132").Value = 4370
$ws.Range("K132").Value = 11437.7145
$ws.Range("L132").Value = 13110
$ws.Range("M132").Value = -8907.7145
$ws.Range("N132").Value = -18170

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 14041.6
$ws.Range("I75").Value = 14041.6
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 14041.6
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -13105.6

$ws.Range("H78").Value = 14041.6
$ws.Range("I78").Value = 14041.6
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 42124.8
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -37444.8

$ws.Range("H80").Value = 689.0833
$ws.Range("I80").Value = 433.54544
$ws.Range("K80").Value = 433.54544
$ws.Range("M80").Value = 564.45456

$ws.Range("H83").Value = 689.0833
$ws.Range("I83").Value = 433.54544
$ws.Range("K83").Value = 2167.7272
$ws.Range("M83").Value = 2824.2728

$ws.Range("H105").Value = 2010
$ws.Range("I105").Value = 2010
$ws.Range("K105").Value = 2010
$ws.Range("M105").Value = -263

$ws.Range("H134").Value = 2288.9375
$ws.Range("I134").Value = 2288.9375
$ws.Range("K134").Value = 6866.8125
$ws.Range("M134").Value = -4331.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1133
$ws.Range("I16").Value = 1133
$ws.Range("K16").Value = 1133
$ws.Range("M16").Value = -846

$ws.Range("H31").Value = 4605.6377
$ws.Range("I31").Value = 3129.25
$ws.Range("J31").Value = 4990.7827
$ws.Range("K31").Value = 3129.25
$ws.Range("L31").Value = 4990.7827
$ws.Range("M31").Value = -2834.25
$ws.Range("N31").Value = -5580.7827

$ws.Range("H34").Value = 4605.6377
$ws.Range("I34").Value = 3129.25
$ws.Range("J34").Value = 4990.7827
$ws.Range("K34").Value = 3129.25
$ws.Range("L34").Value = 4990.7827
$ws.Range("M34").Value = -2927.25
$ws.Range("N34").Value = -5394.7827

$ws.Range("H99").Value = 2438.2
$ws.Range("I99").Value = 2297.75
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2297.75
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -799.75
$ws.Range("N99").Value = -5996

$ws.Range("H113").Value = 1133
$ws.Range("I113").Value = 1133
$ws.Range("K113").Value = 1133
$ws.Range("M113").Value = 1037

$ws.Range("H126").Value = 2438.2
$ws.Range("I126").Value = 2297.75
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6893.25
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4423.25
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 2172.1
$ws.Range("I132").Value = 1840.25
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 5520.75
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -2990.75
$ws.Range("N132").Value = -15558.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 456412.38
$ws.Range("I4").Value = 695297.9
$ws.Range("K4").Value = 2085893.7
$ws.Range("M4").Value = -2085781.7

$ws.Range("H12").Value = 113.882355
$ws.Range("I12").Value = 14.2
$ws.Range("J12").Value = 155.41667
$ws.Range("K12").Value = 42.59999999999999
$ws.Range("L12").Value = 466.25001
$ws.Range("M12").Value = 130.4
$ws.Range("N12").Value = -812.25001

$ws.Range("H18").Value = 4799.9
$ws.Range("I18").Value = 3999.5
$ws.Range("K18").Value = 11998.5
$ws.Range("M18").Value = -11829.5

$ws.Range("H34").Value = 941.05884
$ws.Range("I34").Value = 154.4
$ws.Range("J34").Value = 1268.8334
$ws.Range("K34").Value = 463.2
$ws.Range("L34").Value = 3806.5002
$ws.Range("M34").Value = -379.2
$ws.Range("N34").Value = -3974.5002

$ws.Range("H39").Value = 4854.6665
$ws.Range("J39").Value = 5553.077
$ws.Range("L39").Value = 16659.231
$ws.Range("N39").Value = -17247.231

$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 75000
$ws.Range("N43").Value = -75228

$ws.Range("H70").Value = 5332.3335
$ws.Range("J70").Value = 5332.3335
$ws.Range("L70").Value = 15997.0005
$ws.Range("N70").Value = -16627.0005

$ws.Range("H73").Value = 5332.3335
$ws.Range("J73").Value = 5332.3335
$ws.Range("L73").Value = 15997.0005
$ws.Range("N73").Value = -18181.0005

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0

$ws.Range("H113").Value = 851.44446
$ws.Range("I113").Value = 565.1667
$ws.Range("J113").Value = 1424
$ws.Range("K113").Value = 1695.5001
$ws.Range("L113").Value = 4272
$ws.Range("M113").Value = 474.4999
$ws.Range("N113").Value = -8612

$ws.Range("H119").Value = 2334.5
$ws.Range("I119").Value = 2334.5
$ws.Range("K119").Value = 7003.5
$ws.Range("M119").Value = -2165.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2263.4
$ws.Range("I107").Value = 1161.909
$ws.Range("K107").Value = 1161.909
$ws.Range("M107").Value = 758.0909999999999

$ws.Range("H122").Value = 1300.8572
$ws.Range("I122").Value = 1281.2
$ws.Range("K122").Value = 3843.6
$ws.Range("M122").Value = -1393.6

$ws.Range("H132").Value = 484.57144
$ws.Range("I132").Value = 484.57144
$ws.Range("K132").Value = 1453.71432
$ws.Range("M132").Value = 1076.28568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("N7").Value = 0

$ws.Range("H40").Value = 5797.4
$ws.Range("I40").Value = 5797.4
$ws.Range("K40").Value = 5797.4
$ws.Range("M40").Value = -5661.4

$ws.Range("H46").Value = 1500
$ws.Range("J46").Value = 1500
$ws.Range("L46").Value = 1500
$ws.Range("N46").Value = -1876

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0

$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16996

$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54984

$ws.Range("H103").Value = 37838.555
$ws.Range("J103").Value = 37838.555
$ws.Range("L103").Value = 37838.555
$ws.Range("N103").Value = -40182.555

$ws.Range("H107").Value = 3323.25
$ws.Range("I107").Value = 4197.6665
$ws.Range("K107").Value = 12592.9995
$ws.Range("M107").Value = -10672.9995

$ws.Range("H126").Value = 5088.4585
$ws.Range("I126").Value = 2896.3
$ws.Range("K126").Value = 8688.900000000001
$ws.Range("M126").Value = -6218.900000000001

$ws.Range("H132").Value = 1771.1666
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -9560

